# Restore/update the "From" value for rule R30 on the Rules sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell C10 holds the lower bound ("From") used by rule R30; change it from 18 to 1.
$ws.Range("C10").Value = 1
